# feat: add 2022-Q1 data
#
# The workbook already carries a trailing "总计" (totals) sheet that is
# refreshed every time a new quarter is appended. To add the 2022-Q1
# snapshot we:
#   1. Repurpose the existing "总计" sheet (it becomes the new "2022-Q1"
#      detail sheet) and overwrite its content with the 2022-Q1 fund
#      holdings table.
#   2. Append a brand-new "总计" sheet after it, rebuilt with the
#      2022-Q1 row prepended to the previous totals history.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the "2022-Q1" detail sheet
# ---------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# wipe any previous contents (the old totals table was 4 columns x 5 rows)
$q1.Range("A1:D5").ClearContents()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("B1:H1").Font.Bold = $true
$q1.Range("B1:H1").HorizontalAlignment = -4108
$q1.Range("B1:H1").VerticalAlignment = -4160
$q1.Range("B1:H1").Borders.LineStyle = 1

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "007130"
Set-TextValue $q1.Range("C2") "中庚小盘价值股票"
Set-TextValue $q1.Range("D2") "40.99"
Set-TextValue $q1.Range("E2") "93.10"
Set-TextValue $q1.Range("F2") "2.77"
Set-TextValue $q1.Range("G2") "1.1354"
$q1.Range("H2").Value = 8

$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "007497"
Set-TextValue $q1.Range("C3") "中庚价值灵动灵活配置混合"
Set-TextValue $q1.Range("D3") "24.35"
Set-TextValue $q1.Range("E3") "89.42"
Set-TextValue $q1.Range("F3") "1.96"
Set-TextValue $q1.Range("G3") "0.4773"
$q1.Range("H3").Value = 9

$q1.Range("A4").Value = 2
Set-TextValue $q1.Range("B4") "006323"
Set-TextValue $q1.Range("C4") "合煦智远嘉选混合A"
Set-TextValue $q1.Range("D4") "1.54"
Set-TextValue $q1.Range("E4") "79.45"
Set-TextValue $q1.Range("F4") "3.52"
Set-TextValue $q1.Range("G4") "0.0542"
$q1.Range("H4").Value = 4

$q1.Range("A5").Value = 3
Set-TextValue $q1.Range("B5") "006324"
Set-TextValue $q1.Range("C5") "合煦智远嘉选混合C"
Set-TextValue $q1.Range("D5") "0.59"
Set-TextValue $q1.Range("E5") "79.45"
Set-TextValue $q1.Range("F5") "3.52"
Set-TextValue $q1.Range("G5") "0.0208"
$q1.Range("H5").Value = 4

$q1.Range("A2:A5").Font.Bold = $true
$q1.Range("A2:A5").HorizontalAlignment = -4108
$q1.Range("A2:A5").VerticalAlignment = -4160
$q1.Range("A2:A5").Borders.LineStyle = 1

# ---------------------------------------------------------------
# Step 2: rebuild the "总计" summary sheet after the new 2022-Q1 sheet
# ---------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("B1:D1").Font.Bold = $true
$total.Range("B1:D1").HorizontalAlignment = -4108
$total.Range("B1:D1").VerticalAlignment = -4160
$total.Range("B1:D1").Borders.LineStyle = 1

$rows = @(
    @("2022-Q1", 4, 1.69),
    @("2021-Q4", 3, 2.22),
    @("2021-Q3", 4, 0.28),
    @("2021-Q2", 4, 0.09),
    @("2020-Q4", 1, 0.16)
)

$r = 2
foreach ($row in $rows) {
    $total.Range("A$r").Value = $r - 2
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
    $r = $r + 1
}

$total.Range("A2:A6").Font.Bold = $true
$total.Range("A2:A6").HorizontalAlignment = -4108
$total.Range("A2:A6").VerticalAlignment = -4160
$total.Range("A2:A6").Borders.LineStyle = 1
